$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.23"
$ws.Range("E2").Value = "'1.05%"
$ws.Range("G2").Value = "'23"
$ws.Range("D3").Value = "'36.33"
$ws.Range("E3").Value = "'1.15%"
$ws.Range("G3").Value = "'23"
$ws.Range("D4").Value = "'5.062"
$ws.Range("E4").Value = "'1.30%"
$ws.Range("G4").Value = "'23"
$ws.Range("D5").Value = "'0.08088"
$ws.Range("E5").Value = "'0.28%"
$ws.Range("G5").Value = "'23"
$ws.Range("D6").Value = "'2.021"
$ws.Range("E6").Value = "'6.58%"
$ws.Range("G6").Value = "'23"
$ws.Range("D7").Value = "'7.839"
$ws.Range("E7").Value = "'-0.17%"
$ws.Range("G7").Value = "'23"
$ws.Range("D8").Value = "'0.9280"
$ws.Range("E8").Value = "'-0.67%"
$ws.Range("G8").Value = "'23"
$ws.Range("E9").Value = "'15.07%"
$ws.Range("G9").Value = "'23"
$ws.Range("D10").Value = "'0.1934"
$ws.Range("E10").Value = "'1.81%"
$ws.Range("G10").Value = "'23"
$ws.Range("D11").Value = "'0.09118"
$ws.Range("E11").Value = "'-0.83%"
$ws.Range("G11").Value = "'23"
$ws.Range("D12").Value = "'0.03530"
$ws.Range("E12").Value = "'0.61%"
$ws.Range("G12").Value = "'23"
$ws.Range("D13").Value = "'0.09862"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("G13").Value = "'23"
$ws.Range("D14").Value = "'0.001414"
$ws.Range("E14").Value = "'-0.89%"
$ws.Range("G14").Value = "'23"
$ws.Range("D15").Value = "'0.006284"
$ws.Range("E15").Value = "'-4.04%"
$ws.Range("G15").Value = "'23"
$ws.Range("D16").Value = "'3.849"
$ws.Range("E16").Value = "'6.56%"
$ws.Range("G16").Value = "'23"
$ws.Range("D17").Value = "'4.164"
$ws.Range("E17").Value = "'0.34%"
$ws.Range("G17").Value = "'23"
$ws.Range("G18").Value = "'23"
$ws.Range("D19").Value = "'0.3448"
$ws.Range("G19").Value = "'23"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("G20").Value = "'23"
$ws.Range("D21").Value = "'4.806"
$ws.Range("E21").Value = "'-8.12%"
$ws.Range("G21").Value = "'23"
$ws.Range("D22").Value = "'0.2346"
$ws.Range("E22").Value = "'-7.43%"
$ws.Range("G22").Value = "'23"
$ws.Range("D23").Value = "'0.04372"
$ws.Range("E23").Value = "'-1.06%"
$ws.Range("G23").Value = "'23"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'0.10%"
$ws.Range("G24").Value = "'23"
$ws.Range("E25").Value = "'-11.88%"
$ws.Range("G25").Value = "'23"
$ws.Range("G26").Value = "'23"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("G27").Value = "'23"
$ws.Range("G28").Value = "'23"
$ws.Range("G29").Value = "'23"
$ws.Range("G30").Value = "'23"
$ws.Range("G31").Value = "'23"
$ws.Range("G32").Value = "'23"
$ws.Range("G33").Value = "'23"
$ws.Range("G34").Value = "'23"
$ws.Range("G35").Value = "'23"
$ws.Range("G36").Value = "'23"
$ws.Range("G37").Value = "'23"
$ws.Range("G38").Value = "'23"
$ws.Range("D39").Value = "'0.02061"
$ws.Range("E39").Value = "'5.44%"
$ws.Range("G39").Value = "'23"
$ws.Range("D40").Value = "'0.05098"
$ws.Range("E40").Value = "'-1.41%"
$ws.Range("G40").Value = "'23"
$ws.Range("D41").Value = "'0.007470"
$ws.Range("E41").Value = "'-1.24%"
$ws.Range("G41").Value = "'23"
$ws.Range("D42").Value = "'0.01014"
$ws.Range("E42").Value = "'-0.58%"
$ws.Range("G42").Value = "'23"
$ws.Range("E43").Value = "'-0.54%"
$ws.Range("G43").Value = "'23"
$ws.Range("D44").Value = "'0.002125"
$ws.Range("E44").Value = "'-6.97%"
$ws.Range("G44").Value = "'23"
$ws.Range("D45").Value = "'0.009666"
$ws.Range("E45").Value = "'-9.91%"
$ws.Range("G45").Value = "'23"
$ws.Range("D46").Value = "'0.00006310"
$ws.Range("E46").Value = "'-0.78%"
$ws.Range("G46").Value = "'23"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("G47").Value = "'23"
$ws.Range("E48").Value = "'1.94%"
$ws.Range("G48").Value = "'23"
$ws.Range("G49").Value = "'23"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("G50").Value = "'23"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.21%"
$ws.Range("G51").Value = "'23"
